$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.027700543403625
$ws.Range("B1").Value = 2.068806171417236
$ws.Range("C1").Value = 2.879647493362427
$ws.Range("D1").Value = 3.488574981689453
$ws.Range("E1").Value = 2.0228590965271
